$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "PSQL Performance" sheet between "Tests" and "SQL Parser"
# ---------------------------------------------------------------------------
$wsTests = $wb.Worksheets.Item("Tests")
$wsPerf = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsTests)
$wsPerf.Name = "PSQL Performance"

# ---------------------------------------------------------------------------
# 2. Populate "PSQL Performance" with the TPCH original/decorrelated timings
# ---------------------------------------------------------------------------
$wsPerf.Range("A1").Value = "TPCH"
$wsPerf.Range("B1").Value = "original"
$wsPerf.Range("C1").Value = "decorrelated"

$wsPerf.Range("A2").Value = "Q2"
$wsPerf.Range("B2").Value = " 226,738 ms"
$wsPerf.Range("C2").Value = "147,444 ms"

$wsPerf.Range("A3").Value = "Q4"
$wsPerf.Range("B3").Value = "183,439 ms"
$wsPerf.Range("C3").Value = "186,085 ms"

$wsPerf.Range("A4").Value = "Q17"
$wsPerf.Range("B4").Value = "2412726,898 ms (40:12,727)"
$wsPerf.Range("C4").Value = "2363,854 ms (00:02,364)"

$wsPerf.Range("A5").Value = "Q20"
$wsPerf.Range("B5").Value = "dnf"
$wsPerf.Range("C5").Value = "1467,400 ms (00:01,467)"

$wsPerf.Range("A6").Value = "Q21"
$wsPerf.Range("B6").Value = "498,419 ms"
$wsPerf.Range("C6").Value = "error"

$wsPerf.Range("A7").Value = "Q22"
$wsPerf.Range("B7").Value = "205,010 ms"
$wsPerf.Range("C7").Value = "329501,595 ms (05:29,502)"
$wsPerf.Range("D7").Value = 'Postgres "not exists" better than "not in" due to null check not in https://explainextended.com/2009/09/16/not-in-vs-not-exists-vs-left-join-is-null-postgresql/'

$wsPerf.Columns.Item(2).ColumnWidth = 25
$wsPerf.Columns.Item(3).ColumnWidth = 22

# ---------------------------------------------------------------------------
# 3. "Meeting Notes" sheet: fix a typo and append the new meeting entries
# ---------------------------------------------------------------------------
$wsNotes = $wb.Worksheets.Item("Meeting Notes")

# typo fix: "PostreSQL" -> "PostgreSQL"
$wsNotes.Range("B30").Value = "PostgreSQL instance with TPCH"

$wsNotes.Range("A31").Value = "08.12.2022"
$wsNotes.Range("B31").Value = "Only missing join predicate and decoupling in Neumann algo"
$wsNotes.Range("C31").Value = "Neumann algo, how get predicate in first step?"

$wsNotes.Range("C32").Value = "Confirm: When decoupling, Neumanns selection is not necessary"

$wsNotes.Range("C33").Value = "Q2 needs select distinct, due to non-equi predicates"

$wsNotes.Range("B34").Value = "PostgreSQL performance discussion"
$wsNotes.Range("C34").Value = "careful when transforming exists/in"

$wsNotes.Range("B35").Value = "First focus on Neumann's algo, then consider exists/in discussion"

$wsNotes.Range("B36").Value = "Idea for potential DBIMP project"

$wsNotes.Range("B37").Value = "Sideway information passing"

$wsNotes.Range("A38").Value = "13.12.2022"
$wsNotes.Range("B38").Value = "why umbra doesn't decouple…"
$wsNotes.Range("C38").Value = 'if "or" an uncorrelated predicate?'

$wsNotes.Range("C39").Value = "if +1 an equi correlated predicate?"

$wsNotes.Range("B40").Value = "Sideway information passing -> CTE"
$wsNotes.Range("C40").Value = "now only supports if attribute names are unique in CTE (with schema, easier to fix this)"

$wsNotes.Range("B41").Value = "What to prioritise"
$wsNotes.Range("C41").Value = "add schema"

$wsNotes.Range("C42").Value = "build UI"

$wsNotes.Range("C43").Value = "unnesting 100% (multiple correlations,all subquery operators, ...)"

$wsNotes.Range("C44").Value = "clean up c-style pointers"

$wsNotes.Range("C45").Value = "visualize RA tree"

# ---------------------------------------------------------------------------
# 4. "Tests" sheet: add the original/decorrelated timing columns, relocate
#    the "comment" column content, and append the "code improvements" list
# ---------------------------------------------------------------------------
$wsTestsData = $wb.Worksheets.Item("Tests")

# Header row: C1/D1 become "original"/"decorrelated"; "comment" moves from
# C1(old) -> E1(new, loses its border style which moves to A14 below)
$wsTestsData.Range("C1").Value = "original"
$wsTestsData.Range("D1").Value = "decorrelated"
$wsTestsData.Range("E1").Value = "comment"

# Row 3: C3's old content moves to E3
$wsTestsData.Range("C3").ClearContents()
$wsTestsData.Range("E3").Value = "can also be used by exists/in with complex predicate?"

# Row 4
$wsTestsData.Range("C4").Value = " 226,738 ms"
$wsTestsData.Range("D4").Value = "147,444 ms"

# Row 5
$wsTestsData.Range("C5").Value = "183,439 ms"
$wsTestsData.Range("D5").Value = "186,085 ms"
$wsTestsData.Range("E5").Value = "https://mariadb.com/kb/en/exists-to-in-optimization/"

# Row 6
$wsTestsData.Range("C6").Value = "2412726,898 ms (40:12,727)"
$wsTestsData.Range("D6").Value = "2363,854 ms (00:02,364)"

# Row 7
$wsTestsData.Range("C7").Value = "dnf"
$wsTestsData.Range("D7").Value = "1467,400 ms (00:01,467)"

# Row 8: C8's old content moves to E8
$oldC8 = $wsTestsData.Range("C8").Value2
$wsTestsData.Range("C8").Value = "498,419 ms"
$wsTestsData.Range("D8").Value = "error (explicit join before implicit)"
$wsTestsData.Range("E8").Value = $oldC8

# Row 9
$wsTestsData.Range("C9").Value = "205,010 ms"
$wsTestsData.Range("D9").Value = "329501,595 ms (05:29,502)"

# Old E1/E2/E3 ("code improvements"/"smart pointers"/"throw exceptions")
# relocate to A14/A15/A16; A14 keeps the bottom-border style that used to
# live on E1.
$wsTestsData.Range("E1").Copy() | Out-Null
$wsTestsData.Range("A14").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$wsTestsData.Range("A14").Value = "code improvements"
$wsTestsData.Range("A15").Value = "smart pointers"
$wsTestsData.Range("A16").Value = "throw exceptions"
$wsTestsData.Range("E2").ClearContents()

$excel.ActiveWorkbook.Worksheets.Item("Tests").Activate()
$excel.ActiveWindow.Zoom = 115
$wsTestsData.Range("B10").Select() | Out-Null

# ---------------------------------------------------------------------------
# 5. View/selection cosmetics (best effort)
# ---------------------------------------------------------------------------
$wsNotes.Activate()
$excel.ActiveWindow.Zoom = 130
$wsNotes.Range("C45").Select() | Out-Null

$wsNotes.Activate()

Write-Host "done"
